$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the e-mail address in column A (rows 133-147): the domain was
# missing the ".mx" suffix ("ygtiripitig@liverpool.com" -> "...com.mx") ---
for ($r = 133; $r -le 147; $r++) {
    $ws.Range("A$r").Value = "ygtiripitig@liverpool.com.mx"
}

# --- Add hyperlinks for the corrected e-mail address, mirroring the
# existing "mailto:" hyperlinks already used for the other e-mail rows ---
$ws.Hyperlinks.Add($ws.Range("A133"), "mailto:ygtiripitig@liverpool.com.mx")
$ws.Hyperlinks.Add($ws.Range("A134:A147"), "mailto:ygtiripitig@liverpool.com.mx", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ygtiripitig@liverpool.com.mx")

# Adding the hyperlinks re-applies the built-in "Hyperlink" cell style to the
# first cell touched by each call; restore the original style so the cells
# keep matching the rest of the column.
$ws.Range("A133").Style = $ws.Range("A135").Style
$ws.Range("A134").Style = $ws.Range("A135").Style

# --- Update the view/selection state to match where the workbook was left:
# scrolled so row 123 is at the top, with G146 as the active selection ---
$ws.Range("G146").Select()
